# "create macro for meta data, implement on TRANSPO"
#
# Adds three new metadata rows (share_url, tweet_text, facebook_art) to the
# TRANSPO sheet, directly above the existing "lead_art" row, and also fills
# in a value for "lead_art" (which previously had no value). The existing
# img_1 hyperlink/row is preserved, just shifted down three rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TRANSPO")

$jpgUrl = "http://s3.amazonaws.com/static.texastribune.org/media/images/SH130-4.jpg"

# 0) Stash a copy of the sheet's existing hyperlink look (the style already
#    used by the img_1 row) onto an unused scratch cell so we can reapply it
#    later - Hyperlinks.Add() likes to stamp its own built-in "Hyperlink"
#    style over whatever a cell had before.
$ws.Range("B7").Copy() | Out-Null
$ws.Range("Z1").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# 1) Drop the pre-existing hyperlink(s). The cell keeps its visual
#    (underline/blue) style even after the hyperlink relationship itself is
#    removed, so we can safely re-create the links later, once the rows
#    below have been shifted into their final places.
foreach ($hl in @($ws.Hyperlinks)) {
    $hl.Delete() | Out-Null
}

# 2) Insert three fresh rows right above the current row 6 ("lead_art").
#    This pushes lead_art / img_1 / img_1_caption / img_1_credit / text_1
#    down from rows 6-10 to rows 9-13.
$ws.Rows("6:8").Insert() | Out-Null

# 3) Populate the three new metadata rows.
$ws.Range("A6").Value = "share_url"
$ws.Range("B6").Value = "trib.it/shale-life"

$ws.Range("A7").Value = "tweet_text"
$ws.Range("B7").Value = "Just a test"

$ws.Range("A8").Value = "facebook_art"
$ws.Range("B8").Value = $jpgUrl

# 4) lead_art (now on row 9) previously had no value - give it one.
$ws.Range("B9").Value = $jpgUrl

# 5) The newly-inserted cells don't carry any explicit style yet; match the
#    plain label/value look used throughout the sheet (style of A5/B5).
$ws.Range("A5").Copy() | Out-Null
$ws.Range("A6:A8").PasteSpecial(-4122) | Out-Null

$ws.Range("B5").Copy() | Out-Null
$ws.Range("B6:B8").PasteSpecial(-4122) | Out-Null

# 6) lead_art's B cell (row 9) also needs the plain style at this point; it
#    will be upgraded to the hyperlink look below along with the rest.
$ws.Range("B9").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# 7) Re-create the hyperlinks on their final rows.
$ws.Hyperlinks.Add($ws.Range("B6"), "http://trib.it/shale-life") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B8"), $jpgUrl) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B9"), $jpgUrl) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B10"), $jpgUrl) | Out-Null

# 8) Restore the sheet's original hyperlink style (stashed in step 0) on all
#    four linked cells, overriding the generic style Hyperlinks.Add applied.
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("B6").PasteSpecial(-4122) | Out-Null
$ws.Range("B8").PasteSpecial(-4122) | Out-Null
$ws.Range("B9").PasteSpecial(-4122) | Out-Null
$ws.Range("B10").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# 9) Remove the scratch cell so it doesn't linger in the saved sheet.
$ws.Range("Z1").Clear() | Out-Null

Write-Host "TRANSPO metadata rows added"
